$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "No. of R3 Excel Row's to Execute" value for the
# Verify_All_Buckets_ORG_PHONE test case from 100 to 8.
$ws.Range("D2").Value = "8"

# Move the active selection to E17 (was H19).
$ws.Range("E17").Select()
